$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Boiler bounds: the fourth room-coordinate bound moves from (2,0) to (4,0) ---
$ws.Range("A5").Value = "(4,0)"

# --- Chillers block: Amount 4 -> 3, Width 1.7 -> 2 ---
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 2

# --- Boilers block: now reads Amount/Width/Length from the file (2 / 1 / 1) ---
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1

# --- Header row gets a separator: medium bottom border under "Room Coordinates" ---
# and the row grows a touch taller to fit it (matches the look already used on row 4).
$ws.Range("A1").Borders.Item(9).Weight = -4138
$ws.Rows.Item(1).RowHeight = 15.75

# --- Leave the cursor where the author left it ---
[void]$ws.Range("C8").Select()
